# This edit finalizes the trial list for
# scene_cat_exp_2023.2 / 9_scenecat_memory_bedrooms_1:
#  - the 28 stimulus rows (sheet rows 2-29) are re-shuffled into their
#    final randomized presentation order (trial_block/trial_total keep
#    their sequential 1..28 / 27..54 numbering, but each row now carries
#    the stimulus/condition/norming data that belongs to its new slot)
#  - the catch-trial image is swapped from stimuli/catch_11.jpg to
#    stimuli/catch_26.jpg
# Columns A:S for rows 2:29 are rewritten in one shot via a 2-D array so
# the sheet ends up in exactly the finished, sanity-checked state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0:subject_id | 1:task | 2:block_total | 3:block_scene | 4:trial_block | 5:trial_total | 6:target_cat | 7:category | 8:cond_cat | 9:cond_mem | 10:correct_answer | 11:stimulus | 12:conceptual | 13:perceptual | 14:typicality | 15:n | 16:p_typicality | 17:p_conceptual | 18:p_perceptual
$data = New-Object "object[,]" 28,19

# sheet row 2  (trial_block 1, trial_total 27)
$data[0,0] = 9
$data[0,1] = 'memory'
$data[0,2] = 1
$data[0,3] = 1
$data[0,4] = 1
$data[0,5] = 27
$data[0,6] = 'bedrooms'
$data[0,7] = 'bedrooms'
$data[0,8] = 'target'
$data[0,9] = 'old'
$data[0,10] = 'j'
$data[0,11] = 'stimuli/img_3bxjb.png'
$data[0,12] = 87.28571428571429
$data[0,13] = 72.65714285714286
$data[0,14] = 79.97142857142858
$data[0,15] = 35
$data[0,16] = 10
$data[0,17] = 10
$data[0,18] = 10

# sheet row 3  (trial_block 2, trial_total 28)
$data[1,0] = 9
$data[1,1] = 'memory'
$data[1,2] = 1
$data[1,3] = 1
$data[1,4] = 2
$data[1,5] = 28
$data[1,6] = 'bedrooms'
$data[1,7] = 'bedrooms'
$data[1,8] = 'target'
$data[1,9] = 'old'
$data[1,10] = 'j'
$data[1,11] = 'stimuli/img_9pfbj.png'
$data[1,12] = 91.27272727272727
$data[1,13] = 80.0909090909091
$data[1,14] = 85.68181818181819
$data[1,15] = 33
$data[1,16] = 10
$data[1,17] = 10
$data[1,18] = 10

# sheet row 4  (trial_block 3, trial_total 29)
$data[2,0] = 9
$data[2,1] = 'memory'
$data[2,2] = 1
$data[2,3] = 1
$data[2,4] = 3
$data[2,5] = 29
$data[2,6] = 'bedrooms'
$data[2,7] = 'bedrooms'
$data[2,8] = 'target'
$data[2,9] = 'old'
$data[2,10] = 'j'
$data[2,11] = 'stimuli/img_f4jxo.png'
$data[2,12] = 82.91666666666667
$data[2,13] = 65.52777777777777
$data[2,14] = 74.22222222222223
$data[2,15] = 36
$data[2,16] = 8
$data[2,17] = 8
$data[2,18] = 8

# sheet row 5  (trial_block 4, trial_total 30)
$data[3,0] = 9
$data[3,1] = 'memory'
$data[3,2] = 1
$data[3,3] = 1
$data[3,4] = 4
$data[3,5] = 30
$data[3,6] = 'bedrooms'
$data[3,7] = 'bedrooms'
$data[3,8] = $null
$data[3,9] = 'new'
$data[3,10] = 'f'
$data[3,11] = 'stimuli/img_sltwe.png'
$data[3,12] = 72.025
$data[3,13] = 46.875
$data[3,14] = 59.45
$data[3,15] = 40
$data[3,16] = 5
$data[3,17] = 5
$data[3,18] = 5

# sheet row 6  (trial_block 5, trial_total 31)
$data[4,0] = 9
$data[4,1] = 'memory'
$data[4,2] = 1
$data[4,3] = 1
$data[4,4] = 5
$data[4,5] = 31
$data[4,6] = 'bedrooms'
$data[4,7] = 'bedrooms'
$data[4,8] = 'target'
$data[4,9] = 'old'
$data[4,10] = 'j'
$data[4,11] = 'stimuli/img_2pnl2.png'
$data[4,12] = 6.621621621621622
$data[4,13] = 7.135135135135135
$data[4,14] = 6.878378378378379
$data[4,15] = 37
$data[4,16] = 1
$data[4,17] = 1
$data[4,18] = 1

# sheet row 7  (trial_block 6, trial_total 32)
$data[5,0] = 9
$data[5,1] = 'memory'
$data[5,2] = 1
$data[5,3] = 1
$data[5,4] = 6
$data[5,5] = 32
$data[5,6] = 'bedrooms'
$data[5,7] = 'bedrooms'
$data[5,8] = 'target'
$data[5,9] = 'old'
$data[5,10] = 'j'
$data[5,11] = 'stimuli/img_t4hvr.png'
$data[5,12] = 61.69230769230769
$data[5,13] = 39.76923076923077
$data[5,14] = 50.73076923076923
$data[5,15] = 39
$data[5,16] = 3
$data[5,17] = 3
$data[5,18] = 3

# sheet row 8  (trial_block 7, trial_total 33)
$data[6,0] = 9
$data[6,1] = 'memory'
$data[6,2] = 1
$data[6,3] = 1
$data[6,4] = 7
$data[6,5] = 33
$data[6,6] = 'bedrooms'
$data[6,7] = 'bedrooms'
$data[6,8] = $null
$data[6,9] = 'new'
$data[6,10] = 'f'
$data[6,11] = 'stimuli/img_jge7p.png'
$data[6,12] = 90.42424242424242
$data[6,13] = 75.63636363636364
$data[6,14] = 83.03030303030303
$data[6,15] = 33
$data[6,16] = 10
$data[6,17] = 10
$data[6,18] = 10

# sheet row 9  (trial_block 8, trial_total 34)
$data[7,0] = 9
$data[7,1] = 'memory'
$data[7,2] = 1
$data[7,3] = 1
$data[7,4] = 8
$data[7,5] = 34
$data[7,6] = 'bedrooms'
$data[7,7] = 'bedrooms'
$data[7,8] = $null
$data[7,9] = 'new'
$data[7,10] = 'f'
$data[7,11] = 'stimuli/img_psgf7.png'
$data[7,12] = 26
$data[7,13] = 11.66666666666667
$data[7,14] = 18.83333333333333
$data[7,15] = 36
$data[7,16] = 1
$data[7,17] = 1
$data[7,18] = 1

# sheet row 10  (trial_block 9, trial_total 35)
$data[8,0] = 9
$data[8,1] = 'memory'
$data[8,2] = 1
$data[8,3] = 1
$data[8,4] = 9
$data[8,5] = 35
$data[8,6] = 'bedrooms'
$data[8,7] = 'bedrooms'
$data[8,8] = 'target'
$data[8,9] = 'old'
$data[8,10] = 'j'
$data[8,11] = 'stimuli/img_fnu4h.png'
$data[8,12] = 85.87179487179488
$data[8,13] = 70.71794871794872
$data[8,14] = 78.2948717948718
$data[8,15] = 39
$data[8,16] = 9
$data[8,17] = 9
$data[8,18] = 9

# sheet row 11  (trial_block 10, trial_total 36)
$data[9,0] = 9
$data[9,1] = 'memory'
$data[9,2] = 1
$data[9,3] = 1
$data[9,4] = 10
$data[9,5] = 36
$data[9,6] = 'bedrooms'
$data[9,7] = 'bedrooms'
$data[9,8] = $null
$data[9,9] = 'new'
$data[9,10] = 'f'
$data[9,11] = 'stimuli/img_c4uwt.png'
$data[9,12] = 44.48387096774194
$data[9,13] = 30.06451612903226
$data[9,14] = 37.2741935483871
$data[9,15] = 31
$data[9,16] = 2
$data[9,17] = 2
$data[9,18] = 2

# sheet row 12  (trial_block 11, trial_total 37)
$data[10,0] = 9
$data[10,1] = 'memory'
$data[10,2] = 1
$data[10,3] = 1
$data[10,4] = 11
$data[10,5] = 37
$data[10,6] = 'bedrooms'
$data[10,7] = 'bedrooms'
$data[10,8] = 'target'
$data[10,9] = 'old'
$data[10,10] = 'j'
$data[10,11] = 'stimuli/img_1vq1v.png'
$data[10,12] = 69.42857142857143
$data[10,13] = 46.59523809523809
$data[10,14] = 58.01190476190476
$data[10,15] = 42
$data[10,16] = 5
$data[10,17] = 5
$data[10,18] = 5

# sheet row 13  (trial_block 12, trial_total 38)
$data[11,0] = 9
$data[11,1] = 'memory'
$data[11,2] = 1
$data[11,3] = 1
$data[11,4] = 12
$data[11,5] = 38
$data[11,6] = 'bedrooms'
$data[11,7] = 'bedrooms'
$data[11,8] = 'target'
$data[11,9] = 'old'
$data[11,10] = 'j'
$data[11,11] = 'stimuli/img_ose78.png'
$data[11,12] = 80.19444444444444
$data[11,13] = 60.25
$data[11,14] = 70.22222222222223
$data[11,15] = 36
$data[11,16] = 8
$data[11,17] = 7
$data[11,18] = 7

# sheet row 14  (trial_block 13, trial_total 39)
$data[12,0] = 9
$data[12,1] = 'memory'
$data[12,2] = 1
$data[12,3] = 1
$data[12,4] = 13
$data[12,5] = 39
$data[12,6] = 'bedrooms'
$data[12,7] = 'bedrooms'
$data[12,8] = 'target'
$data[12,9] = 'old'
$data[12,10] = 'j'
$data[12,11] = 'stimuli/img_yteqw.png'
$data[12,12] = 66.83783783783784
$data[12,13] = 43.78378378378378
$data[12,14] = 55.31081081081081
$data[12,15] = 37
$data[12,16] = 4
$data[12,17] = 4
$data[12,18] = 4

# sheet row 15  (trial_block 14, trial_total 40)
$data[13,0] = 9
$data[13,1] = 'memory'
$data[13,2] = 1
$data[13,3] = 1
$data[13,4] = 14
$data[13,5] = 40
$data[13,6] = 'bedrooms'
$data[13,7] = 'bedrooms'
$data[13,8] = 'target'
$data[13,9] = 'old'
$data[13,10] = 'j'
$data[13,11] = 'stimuli/img_jivhq.png'
$data[13,12] = 37
$data[13,13] = 22.26530612244898
$data[13,14] = 29.63265306122449
$data[13,15] = 49
$data[13,16] = 2
$data[13,17] = 2
$data[13,18] = 2

# sheet row 16  (trial_block 15, trial_total 41)
$data[14,0] = 9
$data[14,1] = 'memory'
$data[14,2] = 1
$data[14,3] = 1
$data[14,4] = 15
$data[14,5] = 41
$data[14,6] = 'bedrooms'
$data[14,7] = $null
$data[14,8] = $null
$data[14,9] = 'catch'
$data[14,10] = 'f'
$data[14,11] = 'stimuli/catch_26.jpg'
$data[14,12] = $null
$data[14,13] = $null
$data[14,14] = $null
$data[14,15] = $null
$data[14,16] = $null
$data[14,17] = $null
$data[14,18] = $null

# sheet row 17  (trial_block 16, trial_total 42)
$data[15,0] = 9
$data[15,1] = 'memory'
$data[15,2] = 1
$data[15,3] = 1
$data[15,4] = 16
$data[15,5] = 42
$data[15,6] = 'bedrooms'
$data[15,7] = 'bedrooms'
$data[15,8] = 'target'
$data[15,9] = 'old'
$data[15,10] = 'j'
$data[15,11] = 'stimuli/img_cmyvx.png'
$data[15,12] = 64.25
$data[15,13] = 40.09375
$data[15,14] = 52.171875
$data[15,15] = 32
$data[15,16] = 4
$data[15,17] = 4
$data[15,18] = 4

# sheet row 18  (trial_block 17, trial_total 43)
$data[16,0] = 9
$data[16,1] = 'memory'
$data[16,2] = 1
$data[16,3] = 1
$data[16,4] = 17
$data[16,5] = 43
$data[16,6] = 'bedrooms'
$data[16,7] = 'bedrooms'
$data[16,8] = 'target'
$data[16,9] = 'old'
$data[16,10] = 'j'
$data[16,11] = 'stimuli/img_aweye.png'
$data[16,12] = 53.42105263157895
$data[16,13] = 31.84210526315789
$data[16,14] = 42.63157894736842
$data[16,15] = 38
$data[16,16] = 2
$data[16,17] = 2
$data[16,18] = 2

# sheet row 19  (trial_block 18, trial_total 44)
$data[17,0] = 9
$data[17,1] = 'memory'
$data[17,2] = 1
$data[17,3] = 1
$data[17,4] = 18
$data[17,5] = 44
$data[17,6] = 'bedrooms'
$data[17,7] = 'bedrooms'
$data[17,8] = 'target'
$data[17,9] = 'old'
$data[17,10] = 'j'
$data[17,11] = 'stimuli/img_juob3.png'
$data[17,12] = 79.92105263157895
$data[17,13] = 59.78947368421053
$data[17,14] = 69.85526315789474
$data[17,15] = 38
$data[17,16] = 7
$data[17,17] = 7
$data[17,18] = 7

# sheet row 20  (trial_block 19, trial_total 45)
$data[18,0] = 9
$data[18,1] = 'memory'
$data[18,2] = 1
$data[18,3] = 1
$data[18,4] = 19
$data[18,5] = 45
$data[18,6] = 'bedrooms'
$data[18,7] = 'bedrooms'
$data[18,8] = 'target'
$data[18,9] = 'old'
$data[18,10] = 'j'
$data[18,11] = 'stimuli/img_72fmj.png'
$data[18,12] = 53.87179487179487
$data[18,13] = 36.02564102564103
$data[18,14] = 44.94871794871795
$data[18,15] = 39
$data[18,16] = 3
$data[18,17] = 3
$data[18,18] = 3

# sheet row 21  (trial_block 20, trial_total 46)
$data[19,0] = 9
$data[19,1] = 'memory'
$data[19,2] = 1
$data[19,3] = 1
$data[19,4] = 20
$data[19,5] = 46
$data[19,6] = 'bedrooms'
$data[19,7] = 'bedrooms'
$data[19,8] = 'target'
$data[19,9] = 'old'
$data[19,10] = 'j'
$data[19,11] = 'stimuli/img_gbypq.png'
$data[19,12] = 76.275
$data[19,13] = 51.925
$data[19,14] = 64.1
$data[19,15] = 40
$data[19,16] = 6
$data[19,17] = 6
$data[19,18] = 6

# sheet row 22  (trial_block 21, trial_total 47)
$data[20,0] = 9
$data[20,1] = 'memory'
$data[20,2] = 1
$data[20,3] = 1
$data[20,4] = 21
$data[20,5] = 47
$data[20,6] = 'bedrooms'
$data[20,7] = 'bedrooms'
$data[20,8] = $null
$data[20,9] = 'new'
$data[20,10] = 'f'
$data[20,11] = 'stimuli/img_4wq98.png'
$data[20,12] = 78.48387096774194
$data[20,13] = 58.12903225806452
$data[20,14] = 68.30645161290323
$data[20,15] = 31
$data[20,16] = 7
$data[20,17] = 7
$data[20,18] = 7

# sheet row 23  (trial_block 22, trial_total 48)
$data[21,0] = 9
$data[21,1] = 'memory'
$data[21,2] = 1
$data[21,3] = 1
$data[21,4] = 22
$data[21,5] = 48
$data[21,6] = 'bedrooms'
$data[21,7] = 'bedrooms'
$data[21,8] = 'target'
$data[21,9] = 'old'
$data[21,10] = 'j'
$data[21,11] = 'stimuli/img_z3yzz.png'
$data[21,12] = 71.71052631578948
$data[21,13] = 49.81578947368421
$data[21,14] = 60.76315789473685
$data[21,15] = 38
$data[21,16] = 5
$data[21,17] = 5
$data[21,18] = 5

# sheet row 24  (trial_block 23, trial_total 49)
$data[22,0] = 9
$data[22,1] = 'memory'
$data[22,2] = 1
$data[22,3] = 1
$data[22,4] = 23
$data[22,5] = 49
$data[22,6] = 'bedrooms'
$data[22,7] = 'bedrooms'
$data[22,8] = $null
$data[22,9] = 'new'
$data[22,10] = 'f'
$data[22,11] = 'stimuli/img_zi682.png'
$data[22,12] = 84.6
$data[22,13] = 69.525
$data[22,14] = 77.0625
$data[22,15] = 40
$data[22,16] = 9
$data[22,17] = 9
$data[22,18] = 9

# sheet row 25  (trial_block 24, trial_total 50)
$data[23,0] = 9
$data[23,1] = 'memory'
$data[23,2] = 1
$data[23,3] = 1
$data[23,4] = 24
$data[23,5] = 50
$data[23,6] = 'bedrooms'
$data[23,7] = 'bedrooms'
$data[23,8] = $null
$data[23,9] = 'new'
$data[23,10] = 'f'
$data[23,11] = 'stimuli/img_ozxpp.png'
$data[23,12] = 26.26470588235294
$data[23,13] = 11.47058823529412
$data[23,14] = 18.86764705882353
$data[23,15] = 34
$data[23,16] = 1
$data[23,17] = 1
$data[23,18] = 1

# sheet row 26  (trial_block 25, trial_total 51)
$data[24,0] = 9
$data[24,1] = 'memory'
$data[24,2] = 1
$data[24,3] = 1
$data[24,4] = 25
$data[24,5] = 51
$data[24,6] = 'bedrooms'
$data[24,7] = 'bedrooms'
$data[24,8] = 'target'
$data[24,9] = 'old'
$data[24,10] = 'j'
$data[24,11] = 'stimuli/img_kzg3h.png'
$data[24,12] = 77.02777777777777
$data[24,13] = 56.22222222222222
$data[24,14] = 66.625
$data[24,15] = 36
$data[24,16] = 7
$data[24,17] = 7
$data[24,18] = 7

# sheet row 27  (trial_block 26, trial_total 52)
$data[25,0] = 9
$data[25,1] = 'memory'
$data[25,2] = 1
$data[25,3] = 1
$data[25,4] = 26
$data[25,5] = 52
$data[25,6] = 'bedrooms'
$data[25,7] = 'bedrooms'
$data[25,8] = 'target'
$data[25,9] = 'old'
$data[25,10] = 'j'
$data[25,11] = 'stimuli/img_ic3os.png'
$data[25,12] = 84.79069767441861
$data[25,13] = 66.16279069767442
$data[25,14] = 75.47674418604652
$data[25,15] = 43
$data[25,16] = 9
$data[25,17] = 9
$data[25,18] = 9

# sheet row 28  (trial_block 27, trial_total 53)
$data[26,0] = 9
$data[26,1] = 'memory'
$data[26,2] = 1
$data[26,3] = 1
$data[26,4] = 27
$data[26,5] = 53
$data[26,6] = 'bedrooms'
$data[26,7] = 'bedrooms'
$data[26,8] = 'target'
$data[26,9] = 'old'
$data[26,10] = 'j'
$data[26,11] = 'stimuli/img_anzgh.png'
$data[26,12] = 75.10526315789474
$data[26,13] = 55.76315789473684
$data[26,14] = 65.4342105263158
$data[26,15] = 38
$data[26,16] = 6
$data[26,17] = 6
$data[26,18] = 6

# sheet row 29  (trial_block 28, trial_total 54)
$data[27,0] = 9
$data[27,1] = 'memory'
$data[27,2] = 1
$data[27,3] = 1
$data[27,4] = 28
$data[27,5] = 54
$data[27,6] = 'bedrooms'
$data[27,7] = 'bedrooms'
$data[27,8] = 'target'
$data[27,9] = 'old'
$data[27,10] = 'j'
$data[27,11] = 'stimuli/img_cgdyc.png'
$data[27,12] = 32.93023255813954
$data[27,13] = 14.04651162790698
$data[27,14] = 23.48837209302326
$data[27,15] = 43
$data[27,16] = 1
$data[27,17] = 1
$data[27,18] = 1

$ws.Range("A2:S29").Value = $data

Write-Output "Shuffled bedrooms memory trial list (rows 2:29) and updated catch-trial image."
